$p = $ppt.ActivePresentation

# Insert a new "Title and Content" slide (layout 16 = ppLayoutText / "Title and
# Content") at position 4 - right after "Business Problem" (slide 3) and
# before "Summary" (previously slide 4).
$newSlide = $p.Slides.Add(4, 16)

# --- Title -----------------------------------------------------------
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Business Opportunity"

# --- Body content ------------------------------------------------------
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Build app within real-estate tech company`rZillow, Redfin, `rUsers visit to learn how renovations will increase their home price`rIncreases revenue via product suite`rFunnel more users to platform"

# Paragraph 2 ("Zillow, Redfin, ...") needs "Opendoor" split into its own run
# (it gets flagged by the spell checker in the original deck) followed by the
# rest of the sentence.
$para2 = $body.Paragraphs(2, 1)
$para2.InsertAfter("Opendoor") | Out-Null
$para2.InsertAfter(", Flip, etc.") | Out-Null

# Indent levels: sub-bullets (level 2 / lvl="1") for the app examples line
# and the "funnel more users" line.
$body.Paragraphs(2, 1).IndentLevel = 2
$body.Paragraphs(5, 1).IndentLevel = 2

# Space-after (6pt == spcPts val="600") on every paragraph in the body.
for ($i = 1; $i -le $body.Paragraphs().Count; $i++) {
    $body.Paragraphs($i, 1).ParagraphFormat.SpaceAfter = 6
}
